$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99-138 down to 100-139.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new weekly price record.
$ws.Cells.Item(99, 1).Value = 11
$ws.Cells.Item(99, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(99, 3).Value = "Bíobío"
$ws.Cells.Item(99, 4).Value = 44755
$ws.Cells.Item(99, 5).Value = 8
$ws.Cells.Item(99, 6).Value = 100112032
$ws.Cells.Item(99, 7).Value = "Zapallo italiano"
$ws.Cells.Item(99, 8).Value = "Sin especificar"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 100
$ws.Cells.Item(99, 11).Value = 8000
$ws.Cells.Item(99, 12).Value = 8500
$ws.Cells.Item(99, 13).Value = 8250
$ws.Cells.Item(99, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(99, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(99, 16).Value = 165
$ws.Cells.Item(99, 17).Value = 50
$ws.Cells.Item(99, 18).Value = "Hortaliza"
